# "Case creator needs a date"
# Insert a new "date of collection" question as row 5 of the survey sheet,
# pushing the existing rows (hid, fname, sex, age, ndep, gps_point) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a blank row above the current row 5 ("hid"); Excel copies the
# formatting from the row above (row 4), which is what the target file shows.
$ws.Rows.Item(5).Insert()

# Fill in the new survey question: type | name | label | required
$ws.Cells.Item(5, 1).Value = "date"
$ws.Cells.Item(5, 2).Value = "coll_date"
$ws.Cells.Item(5, 3).Value = "Date of collection"
$ws.Cells.Item(5, 4).Value = "yes"

# Leave the selection on the newly inserted row, mirroring the edit session.
$ws.Rows.Item(5).Select()
